# C5-PowerPoint.pptx edit
#
# 1) The single table on the deck (slide 6) gets a different built-in
#    table style applied (Table Design gallery pick), changing its
#    <a:tableStyleId> from {628C806B-56E4-4D4F-9F87-B2BDF297116C} to
#    {BBC04350-5C3E-4FB7-819F-41BBDA689002}.
#
# 2) The presentation's theme colour scheme is switched from the
#    "Integral" palette back to the stock "Office Theme" palette
#    (all twelve theme colour slots).

$p = $ppt.ActivePresentation

# --- 1) Re-style the table -------------------------------------------------
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle("{BBC04350-5C3E-4FB7-819F-41BBDA689002}")
        }
    }
}

# --- 2) Restore the default "Office Theme" colour scheme -------------------
$design = $p.Designs.Item(1)
$master = $design.SlideMaster
$colors = $master.Theme.ThemeColorScheme

$colors.Colors(1).RGB  = 0         # dk1      000000
$colors.Colors(2).RGB  = 16777215  # lt1      FFFFFF
$colors.Colors(3).RGB  = 6968388   # dk2      44546A
$colors.Colors(4).RGB  = 15132391  # lt2      E7E6E6
$colors.Colors(5).RGB  = 13998939  # accent1  5B9BD5
$colors.Colors(6).RGB  = 3243501   # accent2  ED7D31
$colors.Colors(7).RGB  = 10855845  # accent3  A5A5A5
$colors.Colors(8).RGB  = 49407     # accent4  FFC000
$colors.Colors(9).RGB  = 12874308  # accent5  4472C4
$colors.Colors(10).RGB = 4697456   # accent6  70AD47
$colors.Colors(11).RGB = 12673797  # hlink    0563C1
$colors.Colors(12).RGB = 7491477   # folHlink 954F72
